$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "demo" / "demo_link" columns (K/L) for the last week of slides
$ws.Range("K42").Value = "Demo"
$ws.Range("L42").Value = "28-maps-2"

$ws.Range("K44").Value = "Demo"
$ws.Range("L44").Value = "30-factors-dates"

$ws.Range("K46").Value = "Demo"
$ws.Range("L46").Value = "31-colour"

$ws.Range("K47").Value = "Demo"
$ws.Range("L47").Value = "32-theme"

$ws.Range("K48").Value = "Demo"
$ws.Range("L48").Value = "33-graphics"

# Reflect the user having scrolled down to / selected the bottom of the
# schedule (last week of slides) before saving
$ws.Range("K49").Select()
